# Update CanonCore link to demo site with fake data.
#
# The VEPPLE bullet list gets reworded/merged (copy shuffles down one
# bullet, dropping the old "Architected Events..." paragraph), the Pavers
# Algolia/Lucky Orange bullets get combined into one, and the CanonCore
# project links swap their visible URL text for friendly labels.

$d = $word.ActiveDocument

function Remove-ParagraphByExactText($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

# --- 1. Reword the "Sole front-end and UX developer..." summary bullet ---
$old1 = "Sole front-end and UX developer, collaborating directly on design & product direction, mentoring engineers, and supporting customer success through client onboarding, demos, and feedback sessions. Helped build a virtual experience platform for 30+ UK universities, reaching 1M+ students with 12+ pages and 7 minutes per session."
$new1 = "Sole front-end and UX developer for a virtual experience platform serving 30+ UK universities and 1M+ students, with 12+ pages and 7 minutes per session. Mentored engineers and supported customer success through onboarding, demos, and feedback sessions."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- 2. Drop the old "Architected Events Management System..." bullet
#        entirely (its wording gets reused one bullet further down, see
#        step 3) - remove it *before* that text reappears elsewhere so the
#        exact-text match below stays unambiguous. ---
Remove-ParagraphByExactText $d "Architected Events Management System with Firebase real-time chat, 1-second status engine, and timezone-aware scheduling, delivering 108% increase in views per session." | Out-Null

# --- 3. The "Built production A/B testing framework..." bullet now holds
#        the old "Architected Events Management System..." copy ---
$old3 = "Built production A/B testing framework with GrowthBook SDK, GDPR-compliant tracking, and GTM analytics, enabling data-driven design decisions."
$new3 = "Architected Events Management System with Firebase real-time chat, 1-second status engine, and timezone-aware scheduling, delivering 108% increase in views per session."
$d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# --- 4. The "Developed custom ACF blocks..." bullet absorbs the old
#        A/B testing copy in front of it ---
$old4 = "Developed custom ACF blocks and GraphQL queries for headless WordPress CMS."
$new4 = "Built production A/B testing framework with GrowthBook SDK, GDPR-compliant tracking, and GTM analytics. Developed custom ACF blocks and GraphQL queries for headless WordPress CMS."
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# --- 5. Drop the "Combined Lucky Orange..." bullet *before* folding its
#        wording into the Algolia bullet (step 6), so it can't collide with
#        the freshly-merged text that also contains this phrase. ---
Remove-ParagraphByExactText $d "Combined Lucky Orange with a weather API for location-based recommendations, resulting in 7+% conversion." | Out-Null

# --- 6. Merge the Algolia bullet with the (now-removed) Lucky Orange copy ---
$old6 = "Integrated Algolia search, improving search relevance and driving a 10% increase in search-driven conversions."
$new6 = "Integrated Algolia search driving 10% more search conversions, and combined Lucky Orange with a weather API for location-based recommendations resulting in 7+% conversion."
$d.Content.Find.Execute($old6, $false, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# --- 7/8. CanonCore project links now show friendly labels instead of the
#          raw URLs (addresses/targets are untouched). Go through the
#          Hyperlink object (TextToDisplay) rather than Find/Replace so the
#          existing run formatting (color/underline) on the link text is
#          kept intact instead of being reset. ---
foreach ($h in $d.Hyperlinks) {
    if ($h.Address -eq "http://canoncore.com") {
        $h.TextToDisplay = "Demo Website"
    } elseif ($h.Address -eq "http://github.com/jacobreesgit/CanonCore") {
        $h.TextToDisplay = "Github"
    }
}
